$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'235.58"
$ws.Range("D3").Value = "'22.32"
$ws.Range("D4").Value = "'5.443"
$ws.Range("D5").Value = "'0.05633"
$ws.Range("D6").Value = "'3.380"
$ws.Range("D7").Value = "'6.480"
$ws.Range("D8").Value = "'1.068"
$ws.Range("D9").Value = "'0.7832"
$ws.Range("D10").Value = "'0.1394"
$ws.Range("D11").Value = "'0.07387"
$ws.Range("D13").Value = "'0.02961"
$ws.Range("D15").Value = "'0.001671"
$ws.Range("D16").Value = "'3.251"
$ws.Range("D17").Value = "'0.04736"
$ws.Range("D18").Value = "'0.0005801"
$ws.Range("D19").Value = "'0.006228"
$ws.Range("D20").Value = "'0.005115"
$ws.Range("D21").Value = "'0.001050"
$ws.Range("D23").Value = "'3.911"
$ws.Range("D26").Value = "'0.1055"
$ws.Range("D27").Value = "'0.0004990"
$ws.Range("D40").Value = "'0.04048"
$ws.Range("D41").Value = "'0.006992"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003500"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.009418"
$ws.Range("D45").Value = "'0.00005441"
$ws.Range("D47").Value = "'0.6753"
$ws.Range("D48").Value = "'0.03986"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.01010"
